# feat: add 2022-Q1 data
#
# The workbook has three sheets: "2021-Q3", "2021-Q4", "总计" (a rolling
# summary of holdings-count / holdings-value per quarter).
#
# This change:
#   1. Turns the old "总计" sheet (3rd tab) into a new "2022-Q1" sheet that
#      holds the per-fund holdings detail for 2022-Q1 (same shape as the
#      "2021-Q3"/"2021-Q4" tabs).
#   2. Appends a brand-new "总计" sheet right after it with the same
#      rolling summary as before, plus a new first data row for "2022-Q1".

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Item(2)
$oldTotals = $wb.Worksheets.Item(3)

# Helper: write a value that *looks* numeric into a cell as literal text
# (mirrors typing `'25.59` into Excel), then strip the resulting
# quote-prefix formatting back to Normal so no stray style lingers.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: duplicate the old "总计" sheet *before* touching it, so the
# duplicate inherits its sheet-level formatting (margins, outline props)
# for the new trailing "总计" sheet we'll populate in step 2. Rename the
# original out of the way first so the duplicate can take the "总计" name.
# ---------------------------------------------------------------------
$oldTotals.Name = "2022-Q1-tmp"
$oldTotals.Copy($null, $oldTotals)
$newTotals = $wb.Worksheets.Item(4)
$newTotals.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: populate the new trailing "总计" sheet: same rolling summary
# rows as before, with a new "2022-Q1" row inserted at the top.
# ---------------------------------------------------------------------
$newTotals.Range("A2").Copy()
$newTotals.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$summaryRows = @(
    @("2022-Q1", 4, 4.6),
    @("2021-Q4", 1, 0.5600000000000001),
    @("2021-Q3", 1, 0.44)
)

for ($r = 0; $r -lt $summaryRows.Length; $r++) {
    $row = $summaryRows[$r]
    $excelRow = $r + 2
    $newTotals.Cells.Item($excelRow, 1).Value = $r
    $newTotals.Cells.Item($excelRow, 2).Value = $row[0]
    $newTotals.Cells.Item($excelRow, 3).Value = $row[1]
    $newTotals.Cells.Item($excelRow, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# Step 3: repurpose the original "总计" sheet into the new "2022-Q1"
# detail sheet (same per-fund holdings shape as "2021-Q3"/"2021-Q4").
# ---------------------------------------------------------------------
$oldTotals.Cells.Clear()
$oldTotals.Name = "2022-Q1"
$q1_2022 = $oldTotals

# Reuse the "2021-Q4" sheet's header row / row formatting (fonts, borders,
# alignment) so the new sheet matches its siblings.
$q4.Range("B1:H1").Copy()
$q1_2022.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$q4.Range("A2").Copy()
$q1_2022.Range("A2:A5").PasteSpecial(-4122) | Out-Null

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q1_2022.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("161810", "银华内需精选混合(LOF)", "25.59", "94.71", "7.36", "1.8834", 5),
    @("009394", "银华同力精选混合", "20.03", "94.68", "7.69", "1.5403", 4),
    @("005106", "银华农业产业股票", "13.24", "93.41", "7.15", "0.9467", 4),
    @("180020", "银华成长先锋混合", "3.05", "79.81", "7.55", "0.2303", 2)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q1_2022.Cells.Item($excelRow, 1).Value = $r

    # B: fund code (text), C: fund name (plain string -- no coercion risk),
    # D-G: fund size / stock position / position ratio / held value, all
    # text. H: position rank, a real number.
    Set-TextValue $q1_2022.Cells.Item($excelRow, 2) $row[0]
    $q1_2022.Cells.Item($excelRow, 3).Value = $row[1]
    Set-TextValue $q1_2022.Cells.Item($excelRow, 4) $row[2]
    Set-TextValue $q1_2022.Cells.Item($excelRow, 5) $row[3]
    Set-TextValue $q1_2022.Cells.Item($excelRow, 6) $row[4]
    Set-TextValue $q1_2022.Cells.Item($excelRow, 7) $row[5]
    $q1_2022.Cells.Item($excelRow, 8).Value = $row[6]
}

# Restore the original selection/active sheet.
$q3.Select()
